$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Alpha")
$ws1.Range("A4:J6").EntireRow.Delete()
$ws2 = $wb.Worksheets.Item("Delta")
$ws2.Range("A4:J6").EntireRow.Delete()

# Update the J column strings to combined versions
$ws1.Range("J2").Value = "eyescrunching+jaw"
$ws1.Range("J3").Value = "jaw+raisingeyebrows"
$ws2.Range("J2").Value = "blink+templerun"
$ws2.Range("J3").Value = "blink+sudoku"

# Update numeric values per diff
$ws1.Range("H2").Value = -0.5382516940701341
$ws1.Range("F3").Value = 0.01
$ws1.Range("H3").Value = -0.3122190014324168

$ws2.Range("H2").Value = -2.432000621081636
$ws2.Range("H3").Value = -2.338056636100122
